$wb = $excel.ActiveWorkbook

# --- Update prices on the TimeRelay sheet ---
$ws = $wb.Worksheets.Item("TimeRelay")
$ws.Range("B2").Value = 297
$ws.Range("B3").Value = 297
$ws.Range("B4").Value = 335
$ws.Range("B5").Value = 335

# --- Add a new "Metadata" sheet after "TimeRelay" with availability info ---
$metaSheet = $wb.Worksheets.Add($null, $ws)
$metaSheet.Name = "Metadata"

$metaSheet.Range("A1").Value = "Энергохит"

# Force the date-like text to stay as plain text instead of being parsed as a date
$metaSheet.Range("B1").NumberFormat = "@"
$metaSheet.Range("B1").Value = "07.24.2013"

$metaSheet.Range("C1").Value = " 01.08.2012"
$metaSheet.Range("D1").Value = "Updated prices"

$metaSheet.Range("A1:D1").Select() | Out-Null

# Restore TimeRelay as the active/selected sheet and cell, matching the saved view
$ws.Activate() | Out-Null
$ws.Range("K9").Select() | Out-Null
